$p = $ppt.ActivePresentation
$s3 = $p.Slides.Item(3)
function Dump-Shapes($shapes, $indent) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $txt = ""
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            $txt = $sh.TextFrame.TextRange.Text -replace "`r", "\r" -replace "`n", "\n"
        }
        Write-Host ("{0}{1}: id={2} name='{3}' type={4} text='{5}'" -f $indent, $i, $sh.Id, $sh.Name, $sh.Type, $txt)
        if ($sh.Type -eq 6) {
            Dump-Shapes $sh.GroupItems ($indent + "  ")
        }
    }
}
Dump-Shapes $s3.Shapes ""
